$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = 'Active cases'

# Data rows (cluster name + active cases), alphabetically sorted with "Outbreak" suffix
$ws.Range("A2").Value = '139 Highett St Apartment Complex Richmond Outbreak'
$ws.Range("B2").Value = 12
$ws.Range("A3").Value = '3153 Sacred Heart Community St Kilda Tier 1A Outbreak'
$ws.Range("B3").Value = 21
$ws.Range("A4").Value = '3528 Ottoman Village Aged Care Broadmeadows Outbreak'
$ws.Range("B4").Value = 16
$ws.Range("A5").Value = '3600 Belvedere Aged Care Noble Park Outbreak'
$ws.Range("B5").Value = 32
$ws.Range("A6").Value = '3612 BlueCross Glengowrie Outbreak'
$ws.Range("B6").Value = 51
$ws.Range("A7").Value = '3652 Regis Aged Care Dandenong North Outbreak'
$ws.Range("B7").Value = 15
$ws.Range("A8").Value = '3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak'
$ws.Range("B8").Value = 26
$ws.Range("A9").Value = '3824 Estia Health South Morang Outbreak'
$ws.Range("B9").Value = 43
$ws.Range("A10").Value = '3980 Arcare Keysborough Aged Care Keysborough Outbreak'
$ws.Range("B10").Value = 26
$ws.Range("A11").Value = '4518 Regis Aged Care Fawkner Outbreak'
$ws.Range("B11").Value = 21
$ws.Range("A12").Value = 'ACFS Port Logistics Altona Outbreak'
$ws.Range("B12").Value = 15
$ws.Range("A13").Value = 'Armstrong Creek School Armstrong Creek Outbreak'
$ws.Range("B13").Value = 14
$ws.Range("A14").Value = 'Aspect Autism Spectrum Australia Disability Service Heatherton Outbreak'
$ws.Range("B14").Value = 12
$ws.Range("A15").Value = 'Bespoke Childcare Dingley Village Outbreak'
$ws.Range("B15").Value = 13
$ws.Range("A16").Value = 'Clifton Hill Primary School Clifton Hill Outbreak'
$ws.Range("B16").Value = 12
$ws.Range("A17").Value = 'Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak'
$ws.Range("B17").Value = 23
$ws.Range("A18").Value = 'Elements Childcare Warralily Armstrong Creek Outbreak'
$ws.Range("B18").Value = 22
$ws.Range("A19").Value = 'Honeyeater Hairdressers Bendigo Outbreak'
$ws.Range("B19").Value = 16
$ws.Range("A20").Value = 'Inghams Enterprises Somerville Outbreak'
$ws.Range("B20").Value = 22
$ws.Range("A21").Value = 'JBS Australia Brooklyn Outbreak'
$ws.Range("B21").Value = 12
$ws.Range("A22").Value = 'Launch Housing City Edge Crisis Accommodation South Melbourne Outbreak'
$ws.Range("B22").Value = 10
$ws.Range("A23").Value = 'McQuinns Gym Bendigo Outbreak'
$ws.Range("B23").Value = 20
$ws.Range("A24").Value = 'Metcash Limited Distribution Centre Laverton North Outbreak'
$ws.Range("B24").Value = 18
$ws.Range("A25").Value = 'Metro Tunnel Shine Domain Site Albert Road Construction Site South Melbourne Outbreak'
$ws.Range("B25").Value = 10
$ws.Range("A26").Value = 'Nido Early School Wyndham Vale Outbreak'
$ws.Range("B26").Value = 13
$ws.Range("A27").Value = 'Northern Health The Northern Hospital Epping Outbreak'
$ws.Range("B27").Value = 13
$ws.Range("A28").Value = 'Shawlands Caravan Park Dandenong South Outbreak'
$ws.Range("B28").Value = 15
$ws.Range("A29").Value = 'St Vincents Hospital Emergency Department Melbourne Outbreak'
$ws.Range("B29").Value = 21
$ws.Range("A30").Value = 'TUROSI PTY LTD Thomastown Outbreak'
$ws.Range("B30").Value = 13
$ws.Range("A31").Value = 'The Robin Hood Inn Drouin West Outbreak'
$ws.Range("B31").Value = 27
$ws.Range("A32").Value = 'The Royal Children''s Hospital Parkville Outbreak'
$ws.Range("B32").Value = 10
$ws.Range("A33").Value = 'Werribee Mercy Hospital Emergency Department Outbreak'
$ws.Range("B33").Value = 43
$ws.Range("A34").Value = 'Western Health Sunshine Hospital Emergency Department Outbreak'
$ws.Range("B34").Value = 25
$ws.Range("A35").Value = 'Yarrabah School Aspendale Outbreak'
$ws.Range("B35").Value = 10
$ws.Range("A36").Value = 'Zed Bar Albury Outbreak'
$ws.Range("B36").Value = 12
